$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "schubert-winterreise_51"
$ws.Range("B2").Value = "schubert-winterreise_21"
$ws.Range("C2").Value = "0.233974358974359"
$ws.Range("D2").Value = "[['D#/G', 'A#:7/F', 'D#/G', 'A#:7', 'D#']]"
$ws.Range("E2").Value = "[['F#:maj', 'C#:7/F', 'F#:maj', 'C#:7/F', 'F#:maj']]"
$ws.Range("F2").Value = "[(40.64, 44.84)]"
$ws.Range("G2").Value = "[(38.58, 45.98)]"
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

$ws.Range("A3").Value = "isophonics_53"
$ws.Range("B3").Value = "schubert-winterreise_68"
$ws.Range("C3").Value = "0.3809523809523809"
$ws.Range("D3").Value = "[['A/3', 'E:7', 'A']]"
$ws.Range("E3").Value = "[['G:maj', 'D:7', 'G:maj']]"
$ws.Range("F3").Value = "[(58.557, 62.834)]"
$ws.Range("G3").Value = "[(8.74, 17.66)]"
$ws.Range("H3").ClearContents()
$ws.Range("I3").ClearContents()

$ws.Range("A4").Value = "schubert-winterreise_40"
$ws.Range("B4").Value = "schubert-winterreise_2"
$ws.Range("C4").Value = "0.5397727272727273"
$ws.Range("D4").Value = "[['D:maj', 'A:7', 'D:maj', 'A:7', 'D:maj']]"
$ws.Range("E4").Value = "[['A:maj/E', 'E:7', 'A:maj', 'E:7', 'A:maj']]"
$ws.Range("F4").Value = "[(27.0, 49.86)]"
$ws.Range("G4").Value = "[(20.56, 26.4)]"
$ws.Range("H4").ClearContents()
$ws.Range("I4").Value = "spotify:track:0XfunCHFEeQnzm4NaY8rJr"

$ws.Range("A5").Value = "jaah_27"
$ws.Range("B5").Value = "schubert-winterreise_130"
$ws.Range("C5").Value = "0.2053571428571428"
$ws.Range("D5").Value = "[['A:7', 'D', 'D']]"
$ws.Range("E5").Value = "[['B:7/A', 'E:maj/G#', 'E:maj/B']]"
$ws.Range("F5").Value = "[(9.82, 12.93)]"
$ws.Range("G5").Value = "[(215.66, 222.96)]"
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()

$ws.Range("A6").Value = "schubert-winterreise_0"
$ws.Range("B6").Value = "schubert-winterreise_61"
$ws.Range("C6").Value = "0.06153846153846154"
$ws.Range("D6").Value = "[['B:min', 'E:min/B', 'B:min'], ['B:min/F#', 'F#:7', 'B:min']]"
$ws.Range("E6").Value = "[['G:min/A#', 'C:min', 'G:min/A#'], ['G:min/D', 'D:7', 'G:min']]"
$ws.Range("F6").Value = "[(25.48, 32.58), (74.1, 80.04)]"
$ws.Range("G6").Value = "[(90.6, 101.18), (42.44, 44.36)]"
$ws.Range("H6").Value = "spotify:track:2g41AZ58LFdQLxmWx82ujI"
$ws.Range("I6").Value = "spotify:track:68YORkKP9uvlOQFMZZZwH5"

$ws.Range("A7").Value = "schubert-winterreise_163"
$ws.Range("B7").Value = "isophonics_128"
$ws.Range("C7").Value = "0.2385964912280702"
$ws.Range("D7").Value = "[['F:maj', 'A#:maj', 'F:maj', 'A#:maj']]"
$ws.Range("E7").Value = "[['C', 'F', 'C', 'F']]"
$ws.Range("F7").Value = "[(111.92, 115.66)]"
$ws.Range("G7").Value = "[(10.634761, 19.597664)]"
$ws.Range("H7").Value = "spotify:track:1nvxQGWCnikMK7a4HYQvSx"
$ws.Range("I7").ClearContents()

$ws.Range("A8").Value = "isophonics_200"
$ws.Range("B8").Value = "isophonics_291"
$ws.Range("C8").Value = "0.2728635682158921"
$ws.Range("D8").Value = "[['E', 'A', 'D', 'E', 'A', 'D', 'A']]"
$ws.Range("E8").Value = "[['A', 'D', 'G', 'A', 'D', 'G', 'D']]"
$ws.Range("F8").Value = "[(54.386064, 77.14162)]"
$ws.Range("G8").Value = "[(31.43458, 42.278299)]"
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = "spotify:track:06ypiqmILMdVeaiErMFA91"

$ws.Range("A9").Value = "isophonics_159"
$ws.Range("B9").Value = "schubert-winterreise_157"
$ws.Range("C9").Value = "0.1589673913043478"
$ws.Range("D9").Value = "[['E:7', 'A', 'E/4']]"
$ws.Range("E9").Value = "[['C:7', 'F:maj', 'C:maj']]"
$ws.Range("F9").Value = "[(48.277, 54.481)]"
$ws.Range("G9").Value = "[(40.34, 42.74)]"
$ws.Range("H9").ClearContents()
$ws.Range("I9").Value = "spotify:track:4lrfYSnZmpXdCWuWqVo8L0"

$ws.Range("A10").Value = "isophonics_108"
$ws.Range("B10").Value = "schubert-winterreise_142"
$ws.Range("C10").Value = "0.09646739130434782"
$ws.Range("D10").Value = "[['A', 'E', 'A']]"
$ws.Range("E10").Value = "[['F:maj', 'C:maj', 'F:maj']]"
$ws.Range("F10").Value = "[(19.000294, 29.309954)]"
$ws.Range("G10").Value = "[(41.08, 43.26)]"
$ws.Range("H10").ClearContents()
$ws.Range("I10").Value = "spotify:track:4lrfYSnZmpXdCWuWqVo8L0"

$ws.Range("A11").Value = "isophonics_21"
$ws.Range("B11").Value = "schubert-winterreise_154"
$ws.Range("C11").Value = "0.2015810276679842"
$ws.Range("D11").Value = "[['G:7', 'C', 'C/b7']]"
$ws.Range("E11").Value = "[['E:7', 'A:maj', 'A:maj']]"
$ws.Range("F11").Value = "[(34.041, 38.588)]"
$ws.Range("G11").Value = "[(9.24, 16.18)]"
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = "spotify:track:0XfunCHFEeQnzm4NaY8rJr"

$ws.Range("A12").Value = "schubert-winterreise_157"
$ws.Range("B12").Value = "schubert-winterreise_68"
$ws.Range("C12").Value = "0.2318840579710145"
$ws.Range("D12").Value = "[['F:maj/A', 'C:7', 'F:maj']]"
$ws.Range("E12").Value = "[['G:maj', 'D:7', 'G:maj']]"
$ws.Range("F12").Value = "[(16.2, 18.46)]"
$ws.Range("G12").Value = "[(8.74, 17.66)]"
$ws.Range("H12").Value = "spotify:track:4lrfYSnZmpXdCWuWqVo8L0"
$ws.Range("I12").ClearContents()

$ws.Range("A13").Value = "isophonics_8"
$ws.Range("B13").Value = "isophonics_82"
$ws.Range("C13").Value = "0.1152482269503546"
$ws.Range("D13").Value = "[['A', 'B', 'E']]"
$ws.Range("E13").Value = "[['D', 'E', 'A']]"
$ws.Range("F13").Value = "[(50.046698, 55.352457)]"
$ws.Range("G13").Value = "[(52.750498, 59.809365)]"
$ws.Range("H13").ClearContents()
$ws.Range("I13").Value = "spotify:track:5EzvwjFwdP5Kfl5AZAemzu"

$ws.Range("A14").Value = "schubert-winterreise_72"
$ws.Range("B14").Value = "schubert-winterreise_136"
$ws.Range("C14").Value = "0.1339285714285714"
$ws.Range("D14").Value = "[['C:maj', 'G:7', 'C:maj'], ['C:maj', 'G:maj', 'C:maj']]"
$ws.Range("E14").Value = "[['F:maj/C', 'C:7', 'F:maj'], ['F:maj', 'C:maj', 'F:maj']]"
$ws.Range("F14").Value = "[(0.38, 7.84), (8.34, 10.04)]"
$ws.Range("G14").Value = "[(79.04, 86.54), (2.5, 26.2)]"
$ws.Range("H14").ClearContents()
$ws.Range("I14").ClearContents()

$ws.Range("A15").Value = "isophonics_43"
$ws.Range("B15").Value = "isophonics_218"
$ws.Range("C15").Value = "0.1726973684210526"
$ws.Range("D15").Value = "[['B', 'E', 'A'], ['B:sus4', 'B', 'E']]"
$ws.Range("E15").Value = "[['G', 'C', 'F'], ['G:sus4', 'G', 'C']]"
$ws.Range("F15").Value = "[(46.016712, 58.392947), (26.14043, 31.051451)]"
$ws.Range("G15").Value = "[(5.776, 9.009), (4.928, 8.311)]"
$ws.Range("H15").ClearContents()
$ws.Range("I15").ClearContents()

$ws.Range("A16").Value = "isophonics_242"
$ws.Range("B16").Value = "isophonics_139"
$ws.Range("C16").Value = "0.1306715063520871"
$ws.Range("D16").Value = "[['A', 'B', 'E']]"
$ws.Range("E16").Value = "[['C', 'D', 'G']]"
$ws.Range("F16").Value = "[(44.310045, 50.196303)]"
$ws.Range("G16").Value = "[(51.85331, 57.05458)]"
$ws.Range("H16").Value = "spotify:track:5SUlhldQJtOhUr2GzH5RI7"
$ws.Range("I16").Value = "spotify:track:25yQPHgC35WNnnOUqFhgVR"

$ws.Range("A17").Value = "schubert-winterreise_162"
$ws.Range("B17").Value = "jaah_55"
$ws.Range("C17").Value = "0.1818181818181818"
$ws.Range("D17").Value = "[['D:7', 'G:maj', 'G:maj']]"
$ws.Range("E17").Value = "[['G:7', 'C', 'C']]"
$ws.Range("F17").Value = "[(10.2, 17.38)]"
$ws.Range("G17").Value = "[(47.25, 51.08)]"
$ws.Range("H17").Value = "spotify:track:0XfunCHFEeQnzm4NaY8rJr"
$ws.Range("I17").ClearContents()
